$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing stats (rows 5,6,8,9 -- average vehicular distance columns D,G,J) ---
$ws.Range("D5").Value = 29.2669398500936
$ws.Range("G5").Value = 34.262582825822101
$ws.Range("J5").Value = 33.126599569429402
$ws.Range("D6").Value = 26.3003068216271
$ws.Range("G6").Value = 51.415014296936299
$ws.Range("J6").Value = 55.514835133287697
$ws.Range("D8").Value = 46.070254260528898
$ws.Range("G8").Value = 56.138624095457999
$ws.Range("J8").Value = 14.1038898837209
$ws.Range("D9").Value = 59.897605292479099
$ws.Range("G9").Value = 187.68558105769199
$ws.Range("J9").Value = 70.545204664031601

# --- Update existing stats (rows 10-13 -- average vehicular distance columns D,G,J) ---
$ws.Range("D10").Value = 5813.3583250000001
$ws.Range("G10").Value = 7588.7599821428503
$ws.Range("J10").Value = 5589.2049904761898
$ws.Range("D11").Value = 4800.0823874999996
$ws.Range("G11").Value = 4480.2299326086904
$ws.Range("J11").Value = 5025.5998954545403
$ws.Range("D12").Value = 4560.9952000000003
$ws.Range("G12").Value = 49732.808720000001
$ws.Range("J12").Value = 205924.0526
$ws.Range("D13").Value = 4630.5787454545398
$ws.Range("G13").Value = 8365.4122839999909
$ws.Range("J13").Value = 10547.5343333333

# --- Populate the actuated (second) table, rows 17-26, columns B-K ---
$ws.Range("B17").Value = "590855737#2.97_3"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0

$ws.Range("A18:K18").Style = "Calculation"
$ws.Range("B18").Value = "590855737#2.97_2"
$ws.Range("C18").Value = 503.33333333333297
$ws.Range("D18").Value = 48.092518205980099
$ws.Range("E18").Value = 68.62
$ws.Range("F18").Value = 283.5
$ws.Range("G18").Value = 58.259761429415001
$ws.Range("H18").Value = 57.49
$ws.Range("I18").Value = 253.75
$ws.Range("J18").Value = 49.128383201581002
$ws.Range("K18").Value = 72.48

$ws.Range("A19:K19").Style = "Calculation"
$ws.Range("B19").Value = "590855737#2.97_1"
$ws.Range("C19").Value = 692
$ws.Range("D19").Value = 31.431881529581499
$ws.Range("E19").Value = 72.900000000000006
$ws.Range("F19").Value = 414.666666666666
$ws.Range("G19").Value = 51.5862726071285
$ws.Range("H19").Value = 58.43
$ws.Range("I19").Value = 348
$ws.Range("J19").Value = 49.866908106169198
$ws.Range("K19").Value = 72.3

$ws.Range("B20").Value = "590855737#2.97_0"
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0.33333333333333298
$ws.Range("G20").Value = 133747.63080000001
$ws.Range("H20").Value = 0.03
$ws.Range("I20").Value = 0.25
$ws.Range("J20").Value = 131500.58979999999
$ws.Range("K20").Value = 0.01

$ws.Range("A21:K21").Style = "Calculation"
$ws.Range("B21").Value = "590855742#3.305_0"
$ws.Range("C21").Value = 480.666666666666
$ws.Range("D21").Value = 78.921944206896598
$ws.Range("E21").Value = 22.7
$ws.Range("F21").Value = 455.33333333333297
$ws.Range("G21").Value = 67.203411614317005
$ws.Range("H21").Value = 28.62
$ws.Range("I21").Value = 1044.25
$ws.Range("J21").Value = 34.663753187394697
$ws.Range("K21").Value = 29.19

$ws.Range("A22:K22").Style = "Calculation"
$ws.Range("B22").Value = "590855742#3.305_1"
$ws.Range("C22").Value = 101.666666666666
$ws.Range("D22").Value = 357.77319050847399
$ws.Range("E22").Value = 6.11
$ws.Range("F22").Value = 151.333333333333
$ws.Range("G22").Value = 244.67948303769401
$ws.Range("H22").Value = 10.57
$ws.Range("I22").Value = 138.5
$ws.Range("J22").Value = 393.38181029668402
$ws.Range("K22").Value = 4.29

$ws.Range("B23").Value = "590855742#3.305_2"
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 5825.2674857142802
$ws.Range("E23").Value = 0.4
$ws.Range("F23").Value = 4.8333333333333304
$ws.Range("G23").Value = 7420.9070344827496
$ws.Range("H23").Value = 0.24
$ws.Range("I23").Value = 5.25
$ws.Range("J23").Value = 6655.7124000000003
$ws.Range("K23").Value = 0.19

$ws.Range("B24").Value = "-590855714#1_0"
$ws.Range("C24").Value = 6.3333333333333304
$ws.Range("D24").Value = 5169.1323166666598
$ws.Range("E24").Value = 0.17
$ws.Range("F24").Value = 7.8333333333333304
$ws.Range("G24").Value = 3931.2337255319098
$ws.Range("H24").Value = 0.43
$ws.Range("I24").Value = 5.75
$ws.Range("J24").Value = 5072.3416956521696
$ws.Range("K24").Value = 0.25

$ws.Range("B25").Value = "-590855714#1_1"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 8806.6607499999991
$ws.Range("E25").Value = 0.16
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 55317.168233333301
$ws.Range("H25").Value = 0.1
$ws.Range("I25").Value = 0.25
$ws.Range("J25").Value = 198077.38449999999
$ws.Range("K25").Value = 0

$ws.Range("B26").Value = "8832625#3_0"
$ws.Range("C26").Value = 4.6666666666666599
$ws.Range("D26").Value = 5330.3329538461503
$ws.Range("E26").Value = 0.65
$ws.Range("F26").Value = 4.3333333333333304
$ws.Range("G26").Value = 8415.0853269230702
$ws.Range("H26").Value = 0.55000000000000004
$ws.Range("I26").Value = 3.25
$ws.Range("J26").Value = 9402.5735846153802
$ws.Range("K26").Value = 0.53

# --- Cursor / selection position ---
$ws.Range("J10").Select()
